$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A: narrower width, and give the whole column a left-aligned,
# wrapped-text default look (matches the new narrower "Description" column).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 41.5

# Build the "white row" column-A look (no fill) on A1, and copy it to the
# other white-row description cells (A6, A8, A11) so they all share one xf.
$ws.Range("A1").WrapText = $true
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").Copy() | Out-Null
foreach ($addr in @("A6", "A8", "A11")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Build the "shaded row" column-A look (fill + wrap + left) on A2. Grab the
# existing shaded fill from B2 (same row) first so the fill style is reused
# rather than recreated, then layer wrap/left on top, then copy that look to
# the other shaded-row description cells, including the previously blank
# A3/A4/A5 continuation cells of the "Calcium binds to CaM" group.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("A2").WrapText = $true
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A2").Copy() | Out-Null
foreach ($addr in @("A3", "A4", "A5", "A7", "A10", "A12")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# "CaM binds to CaMKII" is no longer modeled as a bimolecular reaction:
# give rows 6 and 7 explicit rate-constant values instead of leaving the
# Value column blank / referencing elsewhere.
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = 20000
$ws.Range("C6").NumberFormat = "0.00E+00"

$ws.Range("F6").Value = 10000000
$ws.Range("F6").NumberFormat = "0.00E+00"

# C7 sits in a shaded ("customFormat") row, so copy the number format from
# an existing shaded numeric cell (C2) to land on the same scientific +
# shaded style, then set the value.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("C7").Value = 4200000

# ---------------------------------------------------------------------------
# Rows 7 and 8 now contain two-line wrapped descriptions in column A, so
# they need to be taller.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 29

# ---------------------------------------------------------------------------
# Leave the selection where the author left it after making the edit.
# ---------------------------------------------------------------------------
$ws.Range("G8").Select() | Out-Null
